$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Rayon_rotule (mm)"
$ws.Range("B10").Value = 10

$ws.Range("B13").Select()
